$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as text (e.g. thousand-separator
# formatted figures like '93.048.31'). Force each target cell to a text
# number format before assigning, so Excel's automatic type detection
# does not silently coerce values such as "17.80" into the number 17.8.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.048.31"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.450.49"
$ws.Range("E3").Value = "  +4.00%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.39"
$ws.Range("E5").Value = "  +2.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.57"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.42"
$ws.Range("E7").Value = "  +6.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.387"
$ws.Range("E8").Value = "  +3.74%  "

$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.982"
$ws.Range("E10").Value = "  +9.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.446.67"
$ws.Range("E11").Value = "  +3.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.42"
$ws.Range("E12").Value = "  +3.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.198"
$ws.Range("E13").Value = "  +5.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.16"
$ws.Range("E14").Value = "  +5.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.108.30"
$ws.Range("E15").Value = "  +4.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.902.64"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000246"
$ws.Range("E17").Value = "  +2.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.18"
$ws.Range("E18").Value = "  +4.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.462.16"
$ws.Range("E19").Value = "  +4.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.80"
$ws.Range("E20").Value = "  +6.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.86"
$ws.Range("E21").Value = "  +10.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.479"
$ws.Range("E22").Value = "  +10.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.36"
$ws.Range("E23").Value = "  +5.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "499.13"
$ws.Range("E24").Value = "  +3.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.68"
$ws.Range("E25").Value = "  +9.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000180"
$ws.Range("E26").Value = "  +1.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "94.22"
$ws.Range("E27").Value = "  +5.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.05"
$ws.Range("E28").Value = "  +6.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.633.95"
$ws.Range("E29").Value = "  +3.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.89"
$ws.Range("E30").Value = "  +12.18%  "

$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.15"
$ws.Range("E32").Value = "  +1.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.136"
$ws.Range("E33").Value = "  +1.79%  "

$ws.Range("E34").Value = "  +2.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.176"
$ws.Range("E35").Value = "  +5.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.26"
$ws.Range("E36").Value = "  +4.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.547"
$ws.Range("E37").Value = "  +6.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "561.23"
$ws.Range("E38").Value = "  +6.67%  "

$ws.Range("E39").Value = "  +6.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.42"
$ws.Range("E40").Value = "  +2.12%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.916"
$ws.Range("E42").Value = "  +6.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.147"
$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.70"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0412"
$ws.Range("E45").Value = "  +7.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.66"
$ws.Range("E46").Value = "  +1.26%  "

$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.53"
$ws.Range("E47").Value = "  -0.31%  "

$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.36"
$ws.Range("E48").Value = "  +2.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.95"
$ws.Range("E49").Value = "  +1.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.09"
$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.01"
$ws.Range("E51").Value = "  +2.73%  "
